$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.075645208358765
$ws.Range("B1").Value = 3.772316932678223
$ws.Range("C1").Value = 3.345948457717896
$ws.Range("D1").Value = 2.772963285446167
$ws.Range("E1").Value = 1.734041452407837
